# Natmi following Dr Hou advice
# Recompute the Efna2-Epha3 LR-pair table: the previous 3-row result (one row
# per sending cluster, always paired with target cluster "FAPs") is replaced
# by a full 3x3 sending-cluster x target-cluster matrix (rows 2-10), with
# updated statistics for every column.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Efna2"
$ws.Cells.Item(2, 3).Value = "Epha3"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 1.408030333333333
$ws.Cells.Item(2, 8).Value = 4.224091
$ws.Cells.Item(2, 9).Value = 0.3454737251382253
$ws.Cells.Item(2, 10).Value = 0.3454737251382253
$ws.Cells.Item(2, 11).Value = 1
$ws.Cells.Item(2, 12).Value = 0.3333333333333333
$ws.Cells.Item(2, 13).Value = 0.06698166666666668
$ws.Cells.Item(2, 14).Value = 0.200945
$ws.Cells.Item(2, 15).Value = 0.003012576978541733
$ws.Cells.Item(2, 16).Value = 0.003012576978541732
$ws.Cells.Item(2, 17).Value = 0.0943122184438889
$ws.Cells.Item(2, 18).Value = 0.848809965995
$ws.Cells.Item(2, 19).Value = 0.001040766191042472
$ws.Cells.Item(2, 20).Value = 0.001040766191042472
$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Efna2"
$ws.Cells.Item(3, 3).Value = "Epha3"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 1.408030333333333
$ws.Cells.Item(3, 8).Value = 4.224091
$ws.Cells.Item(3, 9).Value = 0.3454737251382253
$ws.Cells.Item(3, 10).Value = 0.3454737251382253
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 22.03620333333333
$ws.Cells.Item(3, 14).Value = 66.10861
$ws.Cells.Item(3, 15).Value = 0.9911034191912899
$ws.Cells.Item(3, 16).Value = 0.9911034191912899
$ws.Cells.Item(3, 17).Value = 31.02764272483444
$ws.Cells.Item(3, 18).Value = 279.24878452351
$ws.Cells.Item(3, 19).Value = 0.3424001902252469
$ws.Cells.Item(3, 20).Value = 0.342400190225247
$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 2).Value = "Efna2"
$ws.Cells.Item(4, 3).Value = "Epha3"
$ws.Cells.Item(4, 4).Value = "sCs"
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 1.408030333333333
$ws.Cells.Item(4, 8).Value = 4.224091
$ws.Cells.Item(4, 9).Value = 0.3454737251382253
$ws.Cells.Item(4, 10).Value = 0.3454737251382253
$ws.Cells.Item(4, 11).Value = 2
$ws.Cells.Item(4, 12).Value = 0.6666666666666666
$ws.Cells.Item(4, 13).Value = 0.130825
$ws.Cells.Item(4, 14).Value = 0.392475
$ws.Cells.Item(4, 15).Value = 0.005884003830168287
$ws.Cells.Item(4, 16).Value = 0.005884003830168287
$ws.Cells.Item(4, 17).Value = 0.1842055683583333
$ws.Cells.Item(4, 18).Value = 1.657850115225
$ws.Cells.Item(4, 19).Value = 0.002032768721935824
$ws.Cells.Item(4, 20).Value = 0.002032768721935824
$ws.Cells.Item(5, 1).Value = "FAPs"
$ws.Cells.Item(5, 2).Value = "Efna2"
$ws.Cells.Item(5, 3).Value = "Epha3"
$ws.Cells.Item(5, 4).Value = "ECs"
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 2.015377
$ws.Cells.Item(5, 8).Value = 6.046131
$ws.Cells.Item(5, 9).Value = 0.494492045565236
$ws.Cells.Item(5, 10).Value = 0.4944920455652361
$ws.Cells.Item(5, 11).Value = 1
$ws.Cells.Item(5, 12).Value = 0.3333333333333333
$ws.Cells.Item(5, 13).Value = 0.06698166666666668
$ws.Cells.Item(5, 14).Value = 0.200945
$ws.Cells.Item(5, 15).Value = 0.003012576978541733
$ws.Cells.Item(5, 16).Value = 0.003012576978541732
$ws.Cells.Item(5, 17).Value = 0.1349933104216667
$ws.Cells.Item(5, 18).Value = 1.214939793795
$ws.Cells.Item(5, 19).Value = 0.00148969535254184
$ws.Cells.Item(5, 20).Value = 0.00148969535254184
$ws.Cells.Item(6, 1).Value = "FAPs"
$ws.Cells.Item(6, 2).Value = "Efna2"
$ws.Cells.Item(6, 3).Value = "Epha3"
$ws.Cells.Item(6, 4).Value = "FAPs"
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 2.015377
$ws.Cells.Item(6, 8).Value = 6.046131
$ws.Cells.Item(6, 9).Value = 0.494492045565236
$ws.Cells.Item(6, 10).Value = 0.4944920455652361
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 22.03620333333333
$ws.Cells.Item(6, 14).Value = 66.10861
$ws.Cells.Item(6, 15).Value = 0.9911034191912899
$ws.Cells.Item(6, 16).Value = 0.9911034191912899
$ws.Cells.Item(6, 17).Value = 44.41125736532333
$ws.Cells.Item(6, 18).Value = 399.70131628791
$ws.Cells.Item(6, 19).Value = 0.4900927571226005
$ws.Cells.Item(6, 20).Value = 0.4900927571226006
$ws.Cells.Item(7, 1).Value = "FAPs"
$ws.Cells.Item(7, 2).Value = "Efna2"
$ws.Cells.Item(7, 3).Value = "Epha3"
$ws.Cells.Item(7, 4).Value = "sCs"
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 2.015377
$ws.Cells.Item(7, 8).Value = 6.046131
$ws.Cells.Item(7, 9).Value = 0.494492045565236
$ws.Cells.Item(7, 10).Value = 0.4944920455652361
$ws.Cells.Item(7, 11).Value = 2
$ws.Cells.Item(7, 12).Value = 0.6666666666666666
$ws.Cells.Item(7, 13).Value = 0.130825
$ws.Cells.Item(7, 14).Value = 0.392475
$ws.Cells.Item(7, 15).Value = 0.005884003830168287
$ws.Cells.Item(7, 16).Value = 0.005884003830168287
$ws.Cells.Item(7, 17).Value = 0.263661696025
$ws.Cells.Item(7, 18).Value = 2.372955264225
$ws.Cells.Item(7, 19).Value = 0.0029095930900936
$ws.Cells.Item(7, 20).Value = 0.0029095930900936
$ws.Cells.Item(8, 1).Value = "sCs"
$ws.Cells.Item(8, 2).Value = "Efna2"
$ws.Cells.Item(8, 3).Value = "Epha3"
$ws.Cells.Item(8, 4).Value = "ECs"
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 0.6522436666666667
$ws.Cells.Item(8, 8).Value = 1.956731
$ws.Cells.Item(8, 9).Value = 0.1600342292965385
$ws.Cells.Item(8, 10).Value = 0.1600342292965385
$ws.Cells.Item(8, 11).Value = 1
$ws.Cells.Item(8, 12).Value = 0.3333333333333333
$ws.Cells.Item(8, 13).Value = 0.06698166666666668
$ws.Cells.Item(8, 14).Value = 0.200945
$ws.Cells.Item(8, 15).Value = 0.003012576978541733
$ws.Cells.Item(8, 16).Value = 0.003012576978541732
$ws.Cells.Item(8, 17).Value = 0.04368836786611112
$ws.Cells.Item(8, 18).Value = 0.393195310795
$ws.Cells.Item(8, 19).Value = 0.0004821154349574209
$ws.Cells.Item(8, 20).Value = 0.0004821154349574209
$ws.Cells.Item(9, 1).Value = "sCs"
$ws.Cells.Item(9, 2).Value = "Efna2"
$ws.Cells.Item(9, 3).Value = "Epha3"
$ws.Cells.Item(9, 4).Value = "FAPs"
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 0.6522436666666667
$ws.Cells.Item(9, 8).Value = 1.956731
$ws.Cells.Item(9, 9).Value = 0.1600342292965385
$ws.Cells.Item(9, 10).Value = 0.1600342292965385
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 22.03620333333333
$ws.Cells.Item(9, 14).Value = 66.10861
$ws.Cells.Item(9, 15).Value = 0.9911034191912899
$ws.Cells.Item(9, 16).Value = 0.9911034191912899
$ws.Cells.Item(9, 17).Value = 14.37297406154556
$ws.Cells.Item(9, 18).Value = 129.35676655391
$ws.Cells.Item(9, 19).Value = 0.1586104718434422
$ws.Cells.Item(9, 20).Value = 0.1586104718434422
$ws.Cells.Item(10, 1).Value = "sCs"
$ws.Cells.Item(10, 2).Value = "Efna2"
$ws.Cells.Item(10, 3).Value = "Epha3"
$ws.Cells.Item(10, 4).Value = "sCs"
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 0.6522436666666667
$ws.Cells.Item(10, 8).Value = 1.956731
$ws.Cells.Item(10, 9).Value = 0.1600342292965385
$ws.Cells.Item(10, 10).Value = 0.1600342292965385
$ws.Cells.Item(10, 11).Value = 2
$ws.Cells.Item(10, 12).Value = 0.6666666666666666
$ws.Cells.Item(10, 13).Value = 0.130825
$ws.Cells.Item(10, 14).Value = 0.392475
$ws.Cells.Item(10, 15).Value = 0.005884003830168287
$ws.Cells.Item(10, 16).Value = 0.005884003830168287
$ws.Cells.Item(10, 17).Value = 0.08532977769166666
$ws.Cells.Item(10, 18).Value = 0.767967999225
$ws.Cells.Item(10, 19).Value = 0.0009416420181388627
$ws.Cells.Item(10, 20).Value = 0.0009416420181388627
